$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D").Insert()

$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$ws.Range("D7").Value = 43465
Write-Host "Done"
